$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 12760.125
$ws.Range("J18").Value = 15856.833
$ws.Range("L18").Value = 15856.833
$ws.Range("N18").Value = -16424.833
$ws.Range("H33").Value = 118.125
$ws.Range("I33").Value = 123
$ws.Range("K33").Value = 123
$ws.Range("M33").Value = 106
$ws.Range("H38").Value = 1507.5
$ws.Range("I38").Value = 1346.6666
$ws.Range("J38").Value = 1990
$ws.Range("K38").Value = 4039.9998
$ws.Range("L38").Value = 5970
$ws.Range("M38").Value = -3667.9998
$ws.Range("N38").Value = -6714
$ws.Range("H98").Value = 1320.6818
$ws.Range("I98").Value = 1335.9524
$ws.Range("K98").Value = 1335.9524
$ws.Range("M98").Value = 162.0476000000001
$ws.Range("H122").Value = 1320.6818
$ws.Range("I122").Value = 1335.9524
$ws.Range("K122").Value = 4007.857199999999
$ws.Range("M122").Value = -1557.857199999999
$ws.Range("H127").Value = 3855.5
$ws.Range("I127").Value = 5033.25
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 15099.75
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = -10139.75
$ws.Range("N127").Value = -14420
$ws.Range("H129").Value = 891.24
$ws.Range("J129").Value = 874.7778
$ws.Range("L129").Value = 2624.3334
$ws.Range("N129").Value = -12624.3334
$ws.Range("H137").Value = 1321.3125
$ws.Range("I137").Value = 1276.0667
$ws.Range("K137").Value = 3828.2001
$ws.Range("M137").Value = -1278.2001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2549.0146
$ws.Range("I32").Value = 1754.2699
$ws.Range("K32").Value = 1754.2699
$ws.Range("M32").Value = -1467.2699
$ws.Range("H122").Value = 4297.3
$ws.Range("I122").Value = 5136.75
$ws.Range("J122").Value = 939.5
$ws.Range("K122").Value = 15410.25
$ws.Range("L122").Value = 2818.5
$ws.Range("M122").Value = -12960.25
$ws.Range("N122").Value = -7718.5
$ws.Range("H132").Value = 1549.5714
$ws.Range("I132").Value = 1019.1111
$ws.Range("J132").Value = 2111.2354
$ws.Range("K132").Value = 3057.3333
$ws.Range("L132").Value = 6333.706200000001
$ws.Range("M132").Value = -527.3332999999998
$ws.Range("N132").Value = -11393.7062

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 60000
$ws.Range("J74").Value = 60000
$ws.Range("L74").Value = 60000
$ws.Range("N74").Value = -61872
$ws.Range("H77").Value = 60000
$ws.Range("J77").Value = 60000
$ws.Range("L77").Value = 180000
$ws.Range("N77").Value = -189360
$ws.Range("H81").Value = 39890
$ws.Range("J81").Value = 39890
$ws.Range("L81").Value = 39890
$ws.Range("N81").Value = -42012
$ws.Range("H84").Value = 39890
$ws.Range("J84").Value = 39890
$ws.Range("L84").Value = 119670
$ws.Range("N84").Value = -130278
$ws.Range("H86").Value = 102020.8
$ws.Range("I86").Value = 1874.4
$ws.Range("J86").Value = 402460
$ws.Range("K86").Value = 1874.4
$ws.Range("L86").Value = 402460
$ws.Range("M86").Value = -751.4000000000001
$ws.Range("N86").Value = -404706
$ws.Range("H89").Value = 102020.8
$ws.Range("I89").Value = 1874.4
$ws.Range("J89").Value = 402460
$ws.Range("K89").Value = 9372
$ws.Range("L89").Value = 2012300
$ws.Range("M89").Value = -3756
$ws.Range("N89").Value = -2023532
$ws.Range("H134").Value = 8514.457
$ws.Range("I134").Value = 8443.4
$ws.Range("J134").Value = 8940.799999999999
$ws.Range("K134").Value = 25330.2
$ws.Range("L134").Value = 26822.4
$ws.Range("M134").Value = -22795.2
$ws.Range("N134").Value = -31892.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 856.6667
$ws.Range("J8").Value = 856.6667
$ws.Range("L8").Value = 856.6667
$ws.Range("N8").Value = -1136.6667
$ws.Range("H23").Value = 65806
$ws.Range("I23").Value = 50000
$ws.Range("J23").Value = 69757.5
$ws.Range("K23").Value = 50000
$ws.Range("L23").Value = 69757.5
$ws.Range("M23").Value = -49760
$ws.Range("N23").Value = -70237.5
$ws.Range("H27").Value = 65806
$ws.Range("I27").Value = 50000
$ws.Range("J27").Value = 69757.5
$ws.Range("K27").Value = 50000
$ws.Range("L27").Value = 69757.5
$ws.Range("M27").Value = -49808
$ws.Range("N27").Value = -70141.5
$ws.Range("H62").Value = 2759.8
$ws.Range("I62").Value = 2574.5
$ws.Range("K62").Value = 2574.5
$ws.Range("M62").Value = -1950.5
$ws.Range("H65").Value = 2759.8
$ws.Range("I65").Value = 2574.5
$ws.Range("K65").Value = 12872.5
$ws.Range("M65").Value = -9752.5
$ws.Range("H70").Value = 20000000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 20000000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H99").Value = 1113879.1
$ws.Range("J99").Value = 3114.125
$ws.Range("L99").Value = 3114.125
$ws.Range("N99").Value = -6110.125
$ws.Range("H105").Value = 1427
$ws.Range("I105").Value = 1329.6666
$ws.Range("K105").Value = 1329.6666
$ws.Range("M105").Value = 417.3334
$ws.Range("H126").Value = 1113879.1
$ws.Range("J126").Value = 3114.125
$ws.Range("L126").Value = 9342.375
$ws.Range("N126").Value = -14282.375
$ws.Range("H132").Value = 2085.5186
$ws.Range("I132").Value = 1219.25
$ws.Range("J132").Value = 3345.5454
$ws.Range("K132").Value = 3657.75
$ws.Range("L132").Value = 10036.6362
$ws.Range("M132").Value = -1127.75
$ws.Range("N132").Value = -15096.6362

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 918.7778
$ws.Range("J5").Value = 921.125
$ws.Range("L5").Value = 2763.375
$ws.Range("N5").Value = -2987.375
$ws.Range("H37").Value = 99999.5
$ws.Range("J37").Value = 99999.5
$ws.Range("L37").Value = 299998.5
$ws.Range("N37").Value = -300222.5
$ws.Range("H87").Value = 25000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H88").Value = 5299.4
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5299.4
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 15898.2
$ws.Range("N88").Value = -16754.2
$ws.Range("H90").Value = 25000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H91").Value = 5299.4
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5299.4
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 15898.2
$ws.Range("N91").Value = -18862.2
$ws.Range("H122").Value = 1008.3
$ws.Range("I122").Value = 932.6667
$ws.Range("J122").Value = 1040.7142
$ws.Range("K122").Value = 8394.0003
$ws.Range("L122").Value = 9366.427799999999
$ws.Range("M122").Value = -5944.0003
$ws.Range("N122").Value = -14266.4278
$ws.Range("H135").Value = 918.7778
$ws.Range("J135").Value = 921.125
$ws.Range("L135").Value = 8290.125
$ws.Range("N135").Value = -13360.125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1618324.9
$ws.Range("I126").Value = 2224841
$ws.Range("K126").Value = 6674523
$ws.Range("M126").Value = -6672053
$ws.Range("H132").Value = 838222.0600000001
$ws.Range("I132").Value = 1013748.94
$ws.Range("K132").Value = 3041246.82
$ws.Range("M132").Value = -3038716.82

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 8493.833000000001
$ws.Range("I20").Value = 4321
$ws.Range("K20").Value = 4321
$ws.Range("M20").Value = -4095
$ws.Range("H40").Value = 4224.45
$ws.Range("I40").Value = 1284
$ws.Range("J40").Value = 9685.286
$ws.Range("K40").Value = 1284
$ws.Range("L40").Value = 9685.286
$ws.Range("M40").Value = -1148
$ws.Range("N40").Value = -9957.286
$ws.Range("H132").Value = 2078.2917
$ws.Range("I132").Value = 2020.1
$ws.Range("K132").Value = 6060.299999999999
$ws.Range("M132").Value = -3530.299999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1921.6389
$ws.Range("I132").Value = 1500.04
$ws.Range("K132").Value = 4500.12
$ws.Range("M132").Value = -1970.12
$ws.Range("H136").Value = 16342145
$ws.Range("I136").Value = 23150144
$ws.Range("K136").Value = 69450432
$ws.Range("M136").Value = -69447882
